$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 18795456
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 18795456
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 56386368
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -56388116
$ws.Range("H72").Value = 18795456
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 18795456
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 169159104
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -169167840
$ws.Range("H137").Value = 14972848
$ws.Range("I137").Value = 23138426
$ws.Range("J137").Value = 2621.6667
$ws.Range("K137").Value = 69415278
$ws.Range("L137").Value = 7865.000100000001
$ws.Range("M137").Value = -69412728
$ws.Range("N137").Value = -12965.0001
$ws.Range("H138").Value = 3551.6191
$ws.Range("I138").Value = 1858.5333
$ws.Range("J138").Value = 3919.6812
$ws.Range("K138").Value = 5575.5999
$ws.Range("L138").Value = 11759.0436
$ws.Range("M138").Value = -435.5999000000002
$ws.Range("N138").Value = -22039.0436
$ws.Range("H141").Value = 1682.5
$ws.Range("I141").Value = 1619
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 4857
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 323
$ws.Range("N141").Value = -16360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18204680
$ws.Range("I32").Value = 32269166
$ws.Range("K32").Value = 32269166
$ws.Range("M32").Value = -32268879
$ws.Range("H74").Value = 2214.0952
$ws.Range("I74").Value = 1142.4286
$ws.Range("J74").Value = 4357.4287
$ws.Range("K74").Value = 1142.4286
$ws.Range("L74").Value = 4357.4287
$ws.Range("M74").Value = -268.4286
$ws.Range("N74").Value = -6105.4287
$ws.Range("H77").Value = 2214.0952
$ws.Range("I77").Value = 1142.4286
$ws.Range("J77").Value = 4357.4287
$ws.Range("K77").Value = 5712.143
$ws.Range("L77").Value = 21787.1435
$ws.Range("M77").Value = -1344.143
$ws.Range("N77").Value = -30523.1435
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 40020
$ws.Range("J132").Value = 40020
$ws.Range("L132").Value = 40020
$ws.Range("N132").Value = -50140
$ws.Range("H134").Value = 2434.56
$ws.Range("I134").Value = 2347.7273
$ws.Range("J134").Value = 3071.3333
$ws.Range("K134").Value = 7043.1819
$ws.Range("L134").Value = 9213.999899999999
$ws.Range("M134").Value = -4508.1819
$ws.Range("N134").Value = -14283.9999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7487.1875
$ws.Range("I132").Value = 13232.429
$ws.Range("K132").Value = 39697.287
$ws.Range("M132").Value = -37167.287
$ws.Range("H134").Value = 2968.5715
$ws.Range("I134").Value = 1320.6364
$ws.Range("J134").Value = 4781.3
$ws.Range("K134").Value = 3961.9092
$ws.Range("L134").Value = 14343.9
$ws.Range("M134").Value = -1426.9092
$ws.Range("N134").Value = -19413.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3135187.8
$ws.Range("I113").Value = 3831746.8
$ws.Range("K113").Value = 11495240.4
$ws.Range("M113").Value = -11493070.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1422.8
$ws.Range("I97").Value = 1438.6471
$ws.Range("J97").Value = 1333
$ws.Range("K97").Value = 1438.6471
$ws.Range("L97").Value = 1333
$ws.Range("M97").Value = -942.6470999999999
$ws.Range("N97").Value = -2325
$ws.Range("H102").Value = 1322.24
$ws.Range("I102").Value = 1275.125
$ws.Range("J102").Value = 1406
$ws.Range("K102").Value = 1275.125
$ws.Range("L102").Value = 1406
$ws.Range("M102").Value = 346.875
$ws.Range("N102").Value = -4650
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1749.8334
$ws.Range("I22").Value = 1499.6666
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1499.6666
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1204.6666
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1749.8334
$ws.Range("I27").Value = 1499.6666
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1499.6666
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1392.6666
$ws.Range("N27").Value = -2214
$ws.Range("H46").Value = 961.25
$ws.Range("I46").Value = 896.6667
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 896.6667
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -708.6667
$ws.Range("N46").Value = -1376
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19463.5
$ws.Range("I62").Value = 22556.2
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 22556.2
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -21932.2
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 19463.5
$ws.Range("I65").Value = 22556.2
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 112781
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -109661
$ws.Range("N65").Value = -26240
$ws.Range("H81").Value = 2913.3333
$ws.Range("I81").Value = 2896
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 5792
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -4731
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 2913.3333
$ws.Range("I84").Value = 2896
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 28960
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -23656
$ws.Range("N84").Value = -40608
$ws.Range("H86").Value = 20000
$ws.Range("J86").Value = 20000
$ws.Range("L86").Value = 20000
$ws.Range("N86").Value = -22246
$ws.Range("H89").Value = 20000
$ws.Range("J89").Value = 20000
$ws.Range("L89").Value = 100000
$ws.Range("N89").Value = -111232
$ws.Range("H132").Value = 2008.196
$ws.Range("I132").Value = 1526.6061
$ws.Range("J132").Value = 2891.111
$ws.Range("K132").Value = 4579.8183
$ws.Range("L132").Value = 8673.332999999999
$ws.Range("M132").Value = -2049.8183
$ws.Range("N132").Value = -13733.333
$ws.Range("H136").Value = 6012.125
$ws.Range("I136").Value = 7518.5
$ws.Range("J136").Value = 2999.375
$ws.Range("K136").Value = 22555.5
$ws.Range("L136").Value = 8998.125
$ws.Range("M136").Value = -20005.5
$ws.Range("N136").Value = -14098.125

Write-Host "Applied all cell updates"